# Add pre-requisites & course aliases to the "IT" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IT")

# Row 10 (SENG2130)                      -> G10: SENG1110 or INFT1004
$ws.Range("G10").Value = "SENG1110 or INFT1004"

# Row 11 (INFT2031)                      -> G11: SENG1110 orINFT1004 or COMP1010
$ws.Range("G11").Value = "SENG1110 orINFT1004 or COMP1010"

# Row 8  (SENG1050 / Web Technologies)   -> H8: COMP1050
$ws.Range("H8").Value = "COMP1050"

# Row 12 (INFT2150)                      -> G12: SENG1050
$ws.Range("G12").Value = "SENG1050"

# Row 13 (SENG2260)                      -> G13: INFT1150 , H13: INFT3150
$ws.Range("G13").Value = "INFT1150"
$ws.Range("H13").Value = "INFT3150"

# Row 14 (INFT3100)                      -> H14: SENG3300 , G14: SENG1050, SENG1110 or INFT1004
$ws.Range("H14").Value = "SENG3300"
$ws.Range("G14").Value = "SENG1050, SENG1110 or INFT1004"

# Row 17 (COMP3851A)                     -> H17: COMP3850
$ws.Range("H17").Value = "COMP3850"

# Row 20 (EBUS3050)                      -> H20: EBUS3010
$ws.Range("H20").Value = "EBUS3010"

# Row 21 (EBUS3030)                      -> G21: COMP1140, SENG1110 or INFT1004
$ws.Range("G21").Value = "COMP1140, SENG1110 or INFT1004"

# Column width adjustments on the IT sheet (widened to fit new content)
$ws.Columns.Item(3).ColumnWidth = 70.83203125
$ws.Columns.Item(7).ColumnWidth = 34
$ws.Columns.Item(8).ColumnWidth = 12.5

# The IT tab (rather than NURSING) is now the active/selected tab, with a
# new active cell/selection.
$ws.Activate()
$ws.Range("G22").Select()
